$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the dbExcel (Neo4j data) and WebExcel (web data) input file names
# from the TC09 companion workbook names to the new TC10 companion workbook names.
$dbName  = "TC10_CCDI_PHS-Accession-phs003111_StudyShrTitle-MolecCharClonal_LibStrat-WGS_Neo4jData.xlsx"
$webName = "TC10_CCDI_PHS-Accession-phs003111_StudyShrTitle-MolecCharClonal_LibStrat-WGS_WebData.xlsx"

$ws.Range("D2:D6").Value = $dbName
$ws.Range("E2:E6").Value = $webName

# Update the FilesTab query (row 6, column B) to limit results to 100 rows instead of 100000.
$filesQuery = $ws.Range("B6").Value2
$filesQuery = $filesQuery.Replace("LIMIT 100000", "LIMIT 100")
$ws.Range("B6").Value = $filesQuery

# Update the active selection shown when the workbook is opened.
$ws.Activate()
$ws.Range("B12").Select()
